# Add the new AbilityConfigChaseHero ability to the enemy abilities list.
# K4 currently holds "Level:ActorConfigs:AbilityConfigHitOnCollision"; update it to
# include the new chase-hero ability as a comma-separated array entry so that it
# matches K5's value, allowing Excel to de-duplicate the shared string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = "Level:ActorConfigs:AbilityConfigHitOnCollision,Level:ActorConfigs:AbilityConfigChaseHero"

$ws.Range("K4").Value = $newValue
$ws.Range("K5").Value = $newValue
